# v2.1 update: refresh elapsed-time metrics & skill ordering, trim stale rows,
# and update the "Checked companies" counter.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Checked companies" counter (T1): 17 -> 12 -------------------
$ws.Range("T1").Value = 12

# --- Update Elapsed Time (K) and Skills (L) values for rows 2-13 ---------
$ws.Range("K2").Value = 18.76
$ws.Range("L2").Value = "manage, job, team, resource, process, capacity, match, experience, candidates, role"

$ws.Range("K3").Value = 31.11
$ws.Range("L3").Value = "data, com, business, work, res, able, team, analysis, skill, skills"

$ws.Range("K4").Value = 64.42
$ws.Range("L4").Value = "data, management, experience, business, work, skills, ability, metrics, portfolio, able"

$ws.Range("K5").Value = 84.11
$ws.Range("L5").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K6").Value = 96.14
$ws.Range("L6").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K7").Value = 107.23
$ws.Range("L7").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K8").Value = 114.88
$ws.Range("L8").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K9").Value = 118.4
$ws.Range("L9").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K10").Value = 121.91
$ws.Range("L10").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K11").Value = 125.67
$ws.Range("L11").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K12").Value = 135.69
$ws.Range("L12").Value = "com, ryder, app, job, work, age, financial, plan, applicant, form"

$ws.Range("K13").Value = 147.77
$ws.Range("L13").Value = "com, work, experience, prime, therapeutics, financial, part, sit, applicant, care"

# --- Remove stale rows 14-18 (Jobgether/Lensa/Lensa/UNFI/ClearCaptions) ---
$ws.Range("A14:T18").EntireRow.Delete()
